# Adds pagination: three new "Elimination" list worksheets after the
# existing "Diabetes Elimination" sheet, each populated the same way as
# the original sheet (single-column list of shared-string values with a
# header-styled first row and body-styled remaining rows).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Reference cells on the existing sheet whose cell styles we reuse:
#   A1  -> header style (s=1)
#   A2  -> plain body style (s=2)
# There is no existing cell with the new "s=4" alternate-border style, so
# we build it once (on the first sheet that needs it) by copying the body
# style and re-coloring its left border, then reuse that cell as the
# donor for every subsequent s=4 cell.

function New-Sheet($afterSheet, $sheetName) {
    $s = $wb.Worksheets.Add($null, $afterSheet)
    $s.Name = $sheetName
    return $s
}

function Set-Row($sheet, $row, $text, $styleDonor) {
    $cell = $sheet.Range("A$row")
    $cell.Value = $text
    $styleDonor.Copy()
    $cell.PasteSpecial(-4122)  # xlPasteFormats
}

# ---------------------------------------------------------------------
# Sheet 2: Hypothyroidism Elimination
# ---------------------------------------------------------------------
$s2 = New-Sheet $ws1 "Hypothyroidism Elimination"

$s2Rows = @(
    "Tofu","Edamame","Tempeh","Cauliflower","Cabbage","Broccoli","Kale",
    "Spinach","Sweet potatoes","Strawberries","Pine nuts","Peanuts",
    "Peaches","Green tea","Coffee","Alcohol","Soy milk","White bread",
    "Cakes, pastries","Fried food","Sugar",
    "Processed food- ham, bacon, salami, sausages","Frozen food","Gluten",
    "Sodas","Energy drinks containing caffeine",
    "Packaged food- noodles, soups, salad dressings, sauces","Candies"
)

for ($i = 0; $i -lt $s2Rows.Count; $i++) {
    $row = $i + 1
    $donor = $(if ($row -eq 1) { $ws1.Range("A1") } else { $ws1.Range("A2") })
    Set-Row $s2 $row $s2Rows[$i] $donor
}

# ---------------------------------------------------------------------
# Sheet 3: Hypertension Elimination
# ---------------------------------------------------------------------
$s3 = New-Sheet $s2 "Hypertension Elimination"

$s3Rows = @(
    "Salty food/snacks(chips,pretzels,crackers)",
    "Caffeine-coffee/tea & many soft drinks",
    "Alcohol",
    "Frozen food, meat (bacon, ham)",
    "Pickles",
    "Processed/canned food",
    "Fried food",
    "Sauces, mayonnaise",
    "Processed meat(bacon,sausages,deli meats)",
    "White rice,",
    "white bread"
)

for ($i = 0; $i -lt $s3Rows.Count; $i++) {
    $row = $i + 1
    $donor = $(if ($row -eq 1) { $ws1.Range("A1") } else { $ws1.Range("A2") })
    Set-Row $s3 $row $s3Rows[$i] $donor
}

# ---------------------------------------------------------------------
# Sheet 4: PCOS Elimination
# ---------------------------------------------------------------------
$s4 = New-Sheet $s3 "PCOS Elimination"

$s4Rows = @(
    "Cakes","Pastries","White bread","Fried food","Pizza","Burger",
    "Carbonated beverages",
    "Sugary foods (sweets, icecreams) and beverages (soda, juices)",
    "Red meat","Processed meat","Dairy","Soy products","Gluten","Pasta",
    "White rice","Doughnuts","Fries","Coffee",
    "Seed oils- vegetable oil, soybean oil, canola oil, rapeseed oil, sunflower oil, safflower oil"
)

# Rows 12 & 13 ("Soy products", "Gluten") use a new alternate-border style
# (cellXfs index 4 / border index 3: like the plain body style but with a
# light-gray left border instead of black). Build that style once by
# copying the plain body style onto A12, then tweak its left border, and
# reuse A12 as the donor for A13.
for ($i = 0; $i -lt $s4Rows.Count; $i++) {
    $row = $i + 1
    if ($row -eq 1) {
        $donor = $ws1.Range("A1")
    } elseif ($row -eq 12) {
        $donor = $ws1.Range("A2")
    } elseif ($row -eq 13) {
        $donor = $s4.Range("A12")
    } else {
        $donor = $ws1.Range("A2")
    }
    Set-Row $s4 $row $s4Rows[$i] $donor
}

# Recolor the left border of the new alternate style (rows 12-13) from
# black to light gray (FFAAAAAA), matching the new border definition.
$s4.Range("A12:A13").Borders.Item(7).Color = 11184810

# Leave the original sheet active/selected, as before the edit.
$ws1.Activate()
